$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: H7 200 -> 0
$ws.Range("H7").Value = 0

# Row 8: D8 200 -> 1000 ; E8 "Paycheck" -> "Cards" ; I8 200 -> 0
$ws.Range("D8").Value = 1000
$ws.Range("I8").Value = 0

# Row 9: D9 1000 -> 800 ; E9 "Paycheck+Tuition" -> "Tuition"
$ws.Range("D9").Value = 800

# Row 10: I10 750 -> 950 (E10 text stays "3rd Paycheck", only sst index shifts)
$ws.Range("I10").Value = 950

# Row 11: D11 900 -> 500 ; E11 "Paycheck+Sick" -> "Sick"
$ws.Range("D11").Value = 500

# Row 13: D13 3600 -> 3400 ; E13 "Paycheck+Taxes" -> "Taxes" ; J13 3500 -> 0
$ws.Range("D13").Value = 3400
$ws.Range("J13").Value = 0

# String edits last, in the order that reproduces the shared-string table order
# observed in the target workbook (Sick, Taxes, Cards, Tuition reuses existing idx)
$ws.Range("E11").Value = "Sick"
$ws.Range("E13").Value = "Taxes"
$ws.Range("E8").Value = "Cards"
$ws.Range("E9").Value = "Tuition"

# Selection change: from L32:L33 (active L33) to D9
$ws.Range("D9").Select() | Out-Null
